$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.745.61"
$ws.Range("E2").Value = "  +5.76%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.262.16"
$ws.Range("E3").Value = "  +4.49%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.20"
$ws.Range("E5").Value = "  +2.87%  "

$ws.Range("E6").Value = "  +3.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "64.32"
$ws.Range("E7").Value = "  +1.06%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.412"
$ws.Range("E9").Value = "  +3.99%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.58"
$ws.Range("E10").Value = "  +2.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0900"
$ws.Range("E11").Value = "  +5.46%  "

$ws.Range("E12").Value = "  +2.36%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.597.33"
$ws.Range("E13").Value = "  +4.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.15"
$ws.Range("E14").Value = "  +0.10%  "

$ws.Range("E15").Value = "  +3.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.828"
$ws.Range("E16").Value = "  +1.66%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.71"
$ws.Range("E17").Value = "  +3.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.253.56"
$ws.Range("E18").Value = "  +4.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.611.78"
$ws.Range("E19").Value = "  +5.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0940"
$ws.Range("E20").Value = "  +10.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.91"
$ws.Range("E21").Value = "  +4.18%  "

$ws.Range("E22").Value = "  -0.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.09"
$ws.Range("E23").Value = "  +9.76%  "

$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.44"
$ws.Range("E25").Value = "  +3.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.36"
$ws.Range("E26").Value = "  +1.50%  "

$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.149"
$ws.Range("E27").Value = "  +5.31%  "

$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.81"
$ws.Range("E28").Value = "  +2.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.55"
$ws.Range("E29").Value = "  -0.32%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.57"
$ws.Range("E30").Value = "  +3.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.46"
$ws.Range("E31").Value = "  +1.97%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.83"
$ws.Range("E32").Value = "  +5.91%  "

$ws.Range("E33").Value = "  +3.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.16"
$ws.Range("E34").Value = "  +9.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.82"
$ws.Range("E35").Value = "  +3.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0639"
$ws.Range("E36").Value = "  +2.67%  "

$ws.Range("E37").Value = "  -2.00%  "

$ws.Range("E38").Value = "  +7.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.47"
$ws.Range("E39").Value = "  +1.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.000264"
$ws.Range("E40").Value = "  +63.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.12"
$ws.Range("E41").Value = "  +19.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.09%  "

$ws.Range("E43").Value = "  +5.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.74"
$ws.Range("E44").Value = "  +12.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.96"
$ws.Range("E45").Value = "  -0.23%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.69"
$ws.Range("E46").Value = "  -0.80%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0988"
$ws.Range("E47").Value = "  +6.66%  "

$ws.Range("E48").Value = "  +1.53%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.509.94"
$ws.Range("E49").Value = "  -0.84%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.13"
$ws.Range("E50").Value = "  +1.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.80"
$ws.Range("E51").Value = "  -0.96%  "
